# Updated symbol list on Mon Jan  2 23:28:46 UTC 2023 with GitHub Actions
#
# Applies the latest price/volume refresh to the crypto ranking sheet.
# Columns D (Price) and E (Volume 1h) hold numeric-looking / percentage
# text values (the sheet stores everything as text), so we force the
# cell's number format to "@" (Text) before assigning the new string -
# otherwise Excel would silently reinterpret "246.13" or "0.74%" as a
# number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Row 2 (BNB) ---
Set-TextValue 2 4 "246.13"
Set-TextValue 2 5 "0.74%"

# --- Row 3 (OKB) ---
Set-TextValue 3 5 "6.25%"

# --- Row 4 (HuobiToken) ---
Set-TextValue 4 4 "5.187"
Set-TextValue 4 5 "1.14%"

# --- Row 5 (Cronos) ---
Set-TextValue 5 4 "0.05735"
Set-TextValue 5 5 "0.97%"

# --- Row 6 (KuCoinToken) ---
Set-TextValue 6 4 "6.569"
Set-TextValue 6 5 "1.13%"

# --- Row 7 (MXToken) ---
Set-TextValue 7 5 "4.80%"

# --- Row 8 (FTXToken) ---
Set-TextValue 8 4 "0.8632"
Set-TextValue 8 5 "0.77%"

# --- Row 9 (WazirX) ---
Set-TextValue 9 4 "0.1364"
Set-TextValue 9 5 "2.33%"

# --- Row 10 (MandalaExchangeToken) ---
Set-TextValue 10 4 "0.07093"
Set-TextValue 10 5 "2.24%"

# --- Row 11 (BitrueCoin) ---
Set-TextValue 11 4 "0.03058"
Set-TextValue 11 5 "7.04%"

# --- Row 12 (BitMartToken) ---
Set-TextValue 12 4 "0.09376"
Set-TextValue 12 5 "-0.22%"

# --- Row 13 (BitForexToken) ---
Set-TextValue 13 4 "0.001535"
Set-TextValue 13 5 "1.05%"

# --- Row 14 (One) ---
Set-TextValue 14 4 "0.0005992"
Set-TextValue 14 5 "-94.12%"

# --- Row 15 (TigerCash) ---
Set-TextValue 15 4 "0.005964"
Set-TextValue 15 5 "-4.01%"

# --- Row 16 (UpBots) ---
Set-TextValue 16 5 "5,226.20%"

# --- Row 17 (LEO) ---
Set-TextValue 17 5 "-0.52%"

# --- Row 18 (GateToken) ---
Set-TextValue 18 4 "3.106"
Set-TextValue 18 5 "3.18%"

# --- Row 19 (BTSEToken) ---
Set-TextValue 19 4 "2.279"
Set-TextValue 19 5 "0.90%"

# --- Row 20 (BitpandaEcosystemToken) ---
Set-TextValue 20 4 "0.3196"
Set-TextValue 20 5 "1.41%"

# --- Row 21 (LiechtensteinCryptoassetsExchange) ---
Set-TextValue 21 4 "0.03293"
Set-TextValue 21 5 "2.19%"

# --- Row 22 (ProBitToken) ---
Set-TextValue 22 4 "0.1300"
Set-TextValue 22 5 "2.09%"

# --- Row 23 (MCDex) ---
Set-TextValue 23 4 "3.471"
Set-TextValue 23 5 "-2.29%"

# --- Row 24 (CoinExToken) ---
Set-TextValue 24 4 "0.04157"
Set-TextValue 24 5 "1.64%"

# --- Row 25 (ZBToken) ---
Set-TextValue 25 5 "0.50%"

# --- Row 26 (BitKan) ---
Set-TextValue 26 5 "1.05%"

# --- Row 27 (HotbitToken) ---
Set-TextValue 27 5 "11.70%"

# --- Row 28 (NitroEx) ---
Set-TextValue 28 5 "2.61%"

# --- Row 40 (IDEX) ---
Set-TextValue 40 4 "0.03753"
Set-TextValue 40 5 "0.89%"

# --- Rows 41-43: ranking reshuffle ---
# Before: 41=BKEXToken, 42=CEJI, 43=KickToken
# After:  41=KickToken, 42=BKEXToken, 43=CEJI
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 41 4 "0.005742"
Set-TextValue 41 5 "-2.97%"

$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1071"
Set-TextValue 42 5 "1.52%"

$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 43 4 "0.002101"
Set-TextValue 43 5 "-8.64%"

# --- Row 44 (LocalTraders) ---
Set-TextValue 44 4 "0.008465"
Set-TextValue 44 5 "-11.78%"

# --- Row 45 (CoinLion) ---
Set-TextValue 45 4 "0.00005283"
Set-TextValue 45 5 "3.45%"

# --- Row 46 (Kangarootoken) ---
Set-TextValue 46 5 "0.06%"

# --- Row 47 (CoinbaseStockToken) ---
Set-TextValue 47 5 "-43.53%"

# --- Row 48 (BOLO) ---
Set-TextValue 48 4 "0.002251"
Set-TextValue 48 5 "-12.37%"

# --- Row 49 (CryptobidCoin) ---
Set-TextValue 49 5 "0.06%"

# --- Row 50 (SpecialPowerGold) ---
Set-TextValue 50 5 "0.06%"
